# Update market-data-derived columns (H, I, J, K, L, M, N) across several
# worksheets, as produced by the scheduled Universalis price-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 822.42255
$ws.Range("J17").Value = 825.98553
$ws.Range("L17").Value = 2477.95659
$ws.Range("N17").Value = -2813.95659
$ws.Range("H86").Value = 1867.5714
$ws.Range("I86").Value = 1740.8
$ws.Range("J86").Value = 2184.5
$ws.Range("K86").Value = 1740.8
$ws.Range("L86").Value = 2184.5
$ws.Range("M86").Value = -617.8
$ws.Range("N86").Value = -4430.5
$ws.Range("H89").Value = 1867.5714
$ws.Range("I89").Value = 1740.8
$ws.Range("J89").Value = 2184.5
$ws.Range("K89").Value = 8704
$ws.Range("L89").Value = 10922.5
$ws.Range("M89").Value = -3088
$ws.Range("N89").Value = -22154.5
$ws.Range("H98").Value = 2006.6086
$ws.Range("I98").Value = 2053.2727
$ws.Range("J98").Value = 980
$ws.Range("K98").Value = 2053.2727
$ws.Range("L98").Value = 980
$ws.Range("M98").Value = -555.2727
$ws.Range("N98").Value = -3976
$ws.Range("H116").Value = 14946.758
$ws.Range("I116").Value = 14206.842
$ws.Range("J116").Value = 15950.929
$ws.Range("K116").Value = 14206.842
$ws.Range("L116").Value = 15950.929
$ws.Range("M116").Value = -10764.842
$ws.Range("N116").Value = -22834.929
$ws.Range("H122").Value = 2006.6086
$ws.Range("I122").Value = 2053.2727
$ws.Range("J122").Value = 980
$ws.Range("K122").Value = 6159.8181
$ws.Range("L122").Value = 2940
$ws.Range("M122").Value = -3709.8181
$ws.Range("N122").Value = -7840
$ws.Range("H132").Value = 2709674.5
$ws.Range("I132").Value = 2872185.8
$ws.Range("K132").Value = 8616557.399999999
$ws.Range("M132").Value = -8614027.399999999
$ws.Range("H138").Value = 52284.3
$ws.Range("I138").Value = 2404.5789
$ws.Range("J138").Value = 999999
$ws.Range("K138").Value = 7213.736699999999
$ws.Range("L138").Value = 2999997
$ws.Range("M138").Value = -2073.736699999999
$ws.Range("N138").Value = -3010277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20330.928
$ws.Range("I32").Value = 22149.36
$ws.Range("K32").Value = 22149.36
$ws.Range("M32").Value = -21862.36
$ws.Range("H45").Value = 4410.4614
$ws.Range("I45").Value = 3319.5
$ws.Range("J45").Value = 6156
$ws.Range("K45").Value = 3319.5
$ws.Range("L45").Value = 6156
$ws.Range("M45").Value = -2942.5
$ws.Range("N45").Value = -6910
$ws.Range("H74").Value = 111535.22
$ws.Range("I74").Value = 118194.86
$ws.Range("K74").Value = 118194.86
$ws.Range("M74").Value = -117320.86
$ws.Range("H77").Value = 111535.22
$ws.Range("I77").Value = 118194.86
$ws.Range("K77").Value = 590974.3
$ws.Range("M77").Value = -586606.3
$ws.Range("H80").Value = 32500.25
$ws.Range("H83").Value = 32500.25
$ws.Range("H88").Value = 4086.2942
$ws.Range("I88").Value = 1155
$ws.Range("J88").Value = 5307.6665
$ws.Range("K88").Value = 1155
$ws.Range("L88").Value = 5307.6665
$ws.Range("M88").Value = -749
$ws.Range("N88").Value = -6119.6665
$ws.Range("H91").Value = 4086.2942
$ws.Range("I91").Value = 1155
$ws.Range("J91").Value = 5307.6665
$ws.Range("K91").Value = 1155
$ws.Range("L91").Value = 5307.6665
$ws.Range("M91").Value = 249
$ws.Range("N91").Value = -8115.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4547802.5
$ws.Range("I31").Value = 7693430.5
$ws.Range("K31").Value = 7693430.5
$ws.Range("M31").Value = -7693135.5
$ws.Range("H34").Value = 4547802.5
$ws.Range("I34").Value = 7693430.5
$ws.Range("K34").Value = 7693430.5
$ws.Range("M34").Value = -7693228.5
$ws.Range("H134").Value = 1309.7894
$ws.Range("I134").Value = 1173.4
$ws.Range("J134").Value = 1821.25
$ws.Range("K134").Value = 3520.2
$ws.Range("L134").Value = 5463.75
$ws.Range("M134").Value = -985.2000000000003
$ws.Range("N134").Value = -10533.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 247
$ws.Range("I2").Value = 177.27272
$ws.Range("K2").Value = 1063.63632
$ws.Range("M2").Value = -950.6363200000001
$ws.Range("H107").Value = 2461.0588
$ws.Range("J107").Value = 899.75
$ws.Range("L107").Value = 2699.25
$ws.Range("N107").Value = -6539.25
$ws.Range("H109").Value = 1460
$ws.Range("I109").Value = 177.5
$ws.Range("K109").Value = 532.5
$ws.Range("M109").Value = 507.5
$ws.Range("H113").Value = 697.73334
$ws.Range("J113").Value = 744.6923
$ws.Range("L113").Value = 2234.0769
$ws.Range("N113").Value = -6574.0769
$ws.Range("H115").Value = 670
$ws.Range("I115").Value = 226.66667
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 680.00001
$ws.Range("L115").Value = 6000
$ws.Range("M115").Value = 494.99999
$ws.Range("N115").Value = -8350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1280.6285
$ws.Range("I97").Value = 850.3333
$ws.Range("J97").Value = 1926.0714
$ws.Range("K97").Value = 850.3333
$ws.Range("L97").Value = 1926.0714
$ws.Range("M97").Value = -354.3333
$ws.Range("N97").Value = -2918.0714
$ws.Range("H107").Value = 474.75
$ws.Range("J107").Value = 799.5
$ws.Range("L107").Value = 799.5
$ws.Range("N107").Value = -4639.5
$ws.Range("H132").Value = 2683.4348
$ws.Range("I132").Value = 2578.1365
$ws.Range("K132").Value = 7734.4095
$ws.Range("M132").Value = -5204.4095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3373333
$ws.Range("J2").Value = 3373333
$ws.Range("L2").Value = 3373333
$ws.Range("N2").Value = -3373557
$ws.Range("H46").Value = 3493.158
$ws.Range("J46").Value = 6596.222
$ws.Range("L46").Value = 6596.222
$ws.Range("N46").Value = -6972.222
$ws.Range("H68").Value = 4548.25
$ws.Range("I68").Value = 3398.1667
$ws.Range("J68").Value = 7998.5
$ws.Range("K68").Value = 3398.1667
$ws.Range("L68").Value = 7998.5
$ws.Range("M68").Value = -2649.1667
$ws.Range("N68").Value = -9496.5
$ws.Range("H71").Value = 4548.25
$ws.Range("I71").Value = 3398.1667
$ws.Range("J71").Value = 7998.5
$ws.Range("K71").Value = 16990.8335
$ws.Range("L71").Value = 39992.5
$ws.Range("M71").Value = -13246.8335
$ws.Range("N71").Value = -47480.5
$ws.Range("H93").Value = 3599.077
$ws.Range("I93").Value = 3838.9
$ws.Range("J93").Value = 2799.6667
$ws.Range("K93").Value = 3838.9
$ws.Range("L93").Value = 2799.6667
$ws.Range("M93").Value = -2590.9
$ws.Range("N93").Value = -5295.6667
$ws.Range("H132").Value = 3172.2856
$ws.Range("I132").Value = 1099
$ws.Range("J132").Value = 4001.6
$ws.Range("K132").Value = 3297
$ws.Range("L132").Value = 12004.8
$ws.Range("M132").Value = -767
$ws.Range("N132").Value = -17064.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 90000
$ws.Range("J123").Value = 90000
$ws.Range("L123").Value = 90000
$ws.Range("N123").Value = -99800
$ws.Range("H126").Value = 1004799.8
$ws.Range("I126").Value = 4499.5
$ws.Range("J126").Value = 1671666.6
$ws.Range("K126").Value = 13498.5
$ws.Range("L126").Value = 5014999.800000001
$ws.Range("M126").Value = -11028.5
$ws.Range("N126").Value = -5019939.800000001
$ws.Range("H132").Value = 2240.6155
$ws.Range("I132").Value = 2222
$ws.Range("J132").Value = 2248.889
$ws.Range("K132").Value = 6666
$ws.Range("L132").Value = 6746.667
$ws.Range("M132").Value = -4136
